$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" header text in F1
$ws.Range("F1").Value = "Last status check on: 23.02.2022 21:45"

# Row 10: swap Cena/Old Cena values, turn Delta Cena into a signed text value,
# and replace the Old Datum serial date with a formatted text timestamp
$ws.Range("B10").Value = 37.7
$ws.Range("C10").Value = 37.4

# "+0.3" looks numeric, so force it in as literal text (quote-prefix) and
# then drop back to the default (unstyled) cell format.
$ws.Range("D10").Value = "'+0.3"
$ws.Range("D10").Style = "Normal"

# E10 previously carried a date-number style; the new value is a plain text
# timestamp, so clear that style back to the default as well.
$ws.Range("E10").Value = "2022-02-23 21:47:11"
$ws.Range("E10").Style = "Normal"
